# Update the cached/displayed regression-table figures in the
# "repeat_loans" worksheet to reflect the refreshed (stabilized) results
# from the external `reg_results/repeat_loans.csv` link.
#
# Each of these cells currently holds a formula such as
# "=[1]repeat_loans!B5" whose last-known (cached) value is what gets
# shown/saved in the cell. The external workbook/CSV backing the link is
# not reachable from this environment, so a genuine link refresh can't be
# performed; instead we push the new, already-known figures into the
# cells via a literal-text formula (e.g. ="0.063") so that:
#   - the cell keeps evaluating to a formula result (type "str"), and
#   - the exact display text (including trailing zeros, asterisks,
#     parentheses, etc.) is preserved verbatim, and
#   - the cell's existing style/number format is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B5"  = "0.063"
    "C5"  = "0.041***"
    "D5"  = "0.024"
    "E5"  = "0.042"
    "F5"  = "0.11**"
    "G5"  = "0.087*"

    "B6"  = "(0.043)"
    "D6"  = "(0.033)"
    "E6"  = "(0.038)"
    "F6"  = "(0.055)"
    "G6"  = "(0.051)"

    "B7"  = "0.050"
    "C7"  = "0.022*"
    "D7"  = "0.028"
    "E7"  = "0.043"
    "F7"  = "0.088*"
    "G7"  = "0.076*"

    "B8"  = "(0.036)"
    "C8"  = "(0.011)"
    "D8"  = "(0.031)"
    "E8"  = "(0.033)"
    "F8"  = "(0.045)"
    "G8"  = "(0.042)"

    "B10" = "6302"
    "C10" = "6302"
    "D10" = "6302"
    "E10" = "6302"
    "F10" = "3032"
    "G10" = "3032"

    "C11" = "0.007"
    "F11" = "0.008"
    "G11" = "0.005"

    "B12" = "0.34"
    "C12" = "0.018"
    "D12" = "0.32"
    "E12" = "0.32"
    "G12" = "0.33"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Formula = '="' + $updates[$addr] + '"'
}
